$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: rename the last two columns
$ws.Range("B1").Value = "Random_value_1"
$ws.Range("C1").Value = "Random_value_2"

# Data rows: column B becomes a numeric "layer count", column C becomes a new numeric value.
# A couple of cells end up blank (empty text) instead of numeric.
$ws.Range("B2").Value = 4
$ws.Range("C2").Value = 345

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 568

$ws.Range("B4").Value = 7
$ws.Range("C4").Value = "'"
$ws.Range("C4").Style = "Normal"

$ws.Range("B5").Value = 7
$ws.Range("C5").Value = 792

$ws.Range("B6").Value = "'"
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").Value = 595

$ws.Range("B7").Value = 7
$ws.Range("C7").Value = 390
